$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1373
$ws.Range("I18").Value = 1373
$ws.Range("K18").Value = 1373
$ws.Range("M18").Value = -1089

$ws.Range("H80").Value = 99.5
$ws.Range("I80").Value = 99.5
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 298.5
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = 699.5
$ws.Range("N80").ClearContents()

$ws.Range("H83").Value = 99.5
$ws.Range("I83").Value = 99.5
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 895.5
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = 4096.5
$ws.Range("N83").ClearContents()

$ws.Range("H98").Value = 195.28572
$ws.Range("I98").Value = 234.90909
$ws.Range("J98").Value = 50
$ws.Range("K98").Value = 234.90909
$ws.Range("L98").Value = 50
$ws.Range("M98").Value = 1263.09091
$ws.Range("N98").Value = -3046

$ws.Range("H112").Value = 2327.4
$ws.Range("J112").Value = 1379
$ws.Range("L112").Value = 4137
$ws.Range("N112").Value = -6353

$ws.Range("H122").Value = 195.28572
$ws.Range("I122").Value = 234.90909
$ws.Range("J122").Value = 50
$ws.Range("K122").Value = 704.72727
$ws.Range("L122").Value = 150
$ws.Range("M122").Value = 1745.27273
$ws.Range("N122").Value = -5050

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 865.9231
$ws.Range("I2").Value = 865.9231
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 865.9231
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -752.9231
$ws.Range("N2").ClearContents()

$ws.Range("H32").Value = 1304.3529
$ws.Range("I32").Value = 1211.6
$ws.Range("K32").Value = 1211.6
$ws.Range("M32").Value = -924.5999999999999
$ws.Range("N32").ClearContents()

$ws.Range("H45").Value = 2339.4119
$ws.Range("I45").Value = 1104.091
$ws.Range("K45").Value = 1104.091
$ws.Range("M45").Value = -727.0909999999999
$ws.Range("N45").ClearContents()

$ws.Range("H88").Value = 1806.8
$ws.Range("I88").Value = 1506
$ws.Range("K88").Value = 1506
$ws.Range("M88").Value = -1100
$ws.Range("N88").ClearContents()

$ws.Range("H91").Value = 1806.8
$ws.Range("I91").Value = 1506
$ws.Range("K91").Value = 1506
$ws.Range("M91").Value = -102
$ws.Range("N91").ClearContents()

$ws.Range("H97").Value = 746.86664
$ws.Range("I97").Value = 728.7857
$ws.Range("K97").Value = 728.7857
$ws.Range("M97").Value = -232.7857
$ws.Range("N97").ClearContents()

$ws.Range("H116").Value = 865.9231
$ws.Range("I116").Value = 865.9231
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 865.9231
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1428.0769
$ws.Range("N116").ClearContents()

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws.Range("H132").Value = 4628.357
$ws.Range("I132").Value = 5489.2
$ws.Range("K132").Value = 16467.6
$ws.Range("M132").Value = -13937.6
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 865.9231
$ws.Range("I3").Value = 865.9231
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 865.9231
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -751.9231
$ws.Range("N3").ClearContents()

$ws.Range("H95").Value = 5989.75
$ws.Range("J95").Value = 5989.75
$ws.Range("L95").Value = 5989.75
$ws.Range("N95").Value = -11481.75

$ws.Range("H105").Value = 1705.909
$ws.Range("I105").Value = 1401.6
$ws.Range("J105").Value = 1959.5
$ws.Range("K105").Value = 1401.6
$ws.Range("L105").Value = 1959.5
$ws.Range("M105").Value = 345.4000000000001
$ws.Range("N105").Value = -5453.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H53").Value = 95995
$ws.Range("J53").Value = 95995
$ws.Range("L53").Value = 95995
$ws.Range("N53").Value = -97209

$ws.Range("H99").Value = 955.5
$ws.Range("I99").Value = 955.5
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 955.5
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = 542.5
$ws.Range("N99").ClearContents()

$ws.Range("H122").Value = 1611.125
$ws.Range("I122").Value = 1611.125
$ws.Range("K122").Value = 4833.375
$ws.Range("M122").Value = -2383.375

$ws.Range("H126").Value = 955.5
$ws.Range("I126").Value = 955.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 2866.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -396.5
$ws.Range("N126").ClearContents()

$ws.Range("H134").Value = 942
$ws.Range("I134").Value = 942
$ws.Range("K134").Value = 2826
$ws.Range("M134").Value = -291

$ws.Range("H135").Value = 70354
$ws.Range("J135").Value = 59999
$ws.Range("L135").Value = 59999
$ws.Range("N135").Value = -70139

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 100.75
$ws.Range("J23").Value = 169.75
$ws.Range("L23").Value = 509.25
$ws.Range("N23").Value = -979.25

$ws.Range("H33").Value = 253.22223
$ws.Range("I33").Value = 300.2857
$ws.Range("K33").Value = 1801.7142
$ws.Range("M33").Value = -1518.7142
$ws.Range("N33").ClearContents()

$ws.Range("H52").Value = 1224
$ws.Range("J52").Value = 1224
$ws.Range("L52").Value = 3672
$ws.Range("N52").Value = -4204

$ws.Range("H107").Value = 321.55554
$ws.Range("I107").Value = 224
$ws.Range("J107").Value = 516.6667
$ws.Range("K107").Value = 672
$ws.Range("L107").Value = 1550.0001
$ws.Range("M107").Value = 1248
$ws.Range("N107").Value = -5390.0001

$ws.Range("H131").Value = 1616.9048
$ws.Range("J131").Value = 2432.6365
$ws.Range("L131").Value = 7297.9095
$ws.Range("N131").Value = -17377.9095

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4296
$ws.Range("J80").Value = 3990
$ws.Range("L80").Value = 3990
$ws.Range("N80").Value = -5986

$ws.Range("H83").Value = 4296
$ws.Range("J83").Value = 3990
$ws.Range("L83").Value = 19950
$ws.Range("N83").Value = -29934

$ws.Range("H102").Value = 1921.4117
$ws.Range("I102").Value = 1711.9333
$ws.Range("K102").Value = 1711.9333
$ws.Range("M102").Value = -89.93329999999992
$ws.Range("N102").ClearContents()

$ws.Range("H113").Value = 2607.4666
$ws.Range("I113").Value = 1176
$ws.Range("K113").Value = 1176
$ws.Range("M113").Value = 994
$ws.Range("N113").ClearContents()

$ws.Range("H122").Value = 626.3333
$ws.Range("I122").Value = 594.5
$ws.Range("K122").Value = 1783.5
$ws.Range("M122").Value = 666.5
$ws.Range("N122").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2998.3333
$ws.Range("I122").Value = 2998.3333
$ws.Range("K122").Value = 8994.999899999999
$ws.Range("M122").Value = -6544.999899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 28944
$ws.Range("J125").Value = 28944
$ws.Range("L125").Value = 28944
$ws.Range("N125").Value = -38784

$ws.Range("H126").Value = 4229.6665
$ws.Range("I126").Value = 2719.75
$ws.Range("J126").Value = 6552.615
$ws.Range("K126").Value = 8159.25
$ws.Range("L126").Value = 19657.845
$ws.Range("M126").Value = -5689.25
$ws.Range("N126").Value = -24597.845
